$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = "-"

# Row 4
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "[-, 'MCT-3A-Elementos de máquinas', -, -]"

# Row 6
$ws.Range("B6").Value = "MCT-2A-M.T.R.M."
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "[-, 'MCT-3A-Elementos de máquinas', -, -]"
$ws.Range("F6").Value = "-"

# Row 7
$ws.Range("B7").Value = "MCT-2A-M.T.R.M."
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "[-, 'MCT-3A-Elementos de máquinas', -, -]"
$ws.Range("F7").Value = "-"

# Row 8
$ws.Range("E8").Value = "[-, 'MCT-3A-Elementos de máquinas', -, -]"
